# Correção nos dados: ajustar cabeçalhos da linha 2 (B2 e F2) para "total",
# alinhando com a correção de dados do início da análise PNAD 2009.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "total"
$ws.Range("F2").Value = "total"
